$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the typo "jeck" -> "jack" in B7 (the one real content edit
# behind the shared-string churn in the diff: "jeck" is dropped from the
# table and "jack" is appended, shifting every later shared-string index
# down by one).
$ws.Range("B7").Value = "jack"

# Selection moved to B7 (active cell), matching the saved sheetView.
$ws.Range("B7").Select()
